# Generate Report for Handoff
#
# A new handoff run has completed: the handoff markdown file and its
# translation artifacts were regenerated under a new GUID / content hash,
# and the handoff timestamps were refreshed. Update the Overview, zh-cn,
# and de-de report sheets (cell values + the matching hyperlink display
# text) to reflect the new handoff.

$wb = $excel.ActiveWorkbook

$oldGuid = "ce0058dc-2f35-40e8-970e-0a9db6247290"
$newGuid = "17295b83-c8e9-45e2-9fe8-940e2078275f"

$oldHash = "5f45226b14e18e274346f4f8ad6905a0a0e3dc43"
$newHash = "3bfc74e2fbbfd210bd704db29f71ddc5a51c0d2d"

$newMdName = "$newGuid.md"
$newZhXlfName = "$newGuid.$newHash.zh-cn.xlf"
$newDeXlfName = "$newGuid.$newHash.de-de.xlf"

$newOverviewDate = "2016-03-24 03:08:59"
$newZhCnDate = "2016-03-24 03:08:55"
$newDeDeDate = "2016-03-24 03:08:59"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# --- Overview sheet ---
$wsOverview.Range("A2").Value = $newMdName
$wsOverview.Range("D2").Value = $newOverviewDate

foreach ($hl in $wsOverview.Hyperlinks) {
    $hl.TextToDisplay = $newMdName
}

# --- zh-cn sheet ---
$wsZhCn.Range("A2").Value = $newMdName
$wsZhCn.Range("D2").Value = $newZhXlfName
$wsZhCn.Range("E2").Value = $newZhCnDate

foreach ($hl in $wsZhCn.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') {
        $hl.TextToDisplay = $newMdName
    } elseif ($addr -eq '$D$2') {
        $hl.TextToDisplay = $newZhXlfName
    }
}

# --- de-de sheet ---
$wsDeDe.Range("A2").Value = $newMdName
$wsDeDe.Range("D2").Value = $newDeXlfName
$wsDeDe.Range("E2").Value = $newDeDeDate

foreach ($hl in $wsDeDe.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') {
        $hl.TextToDisplay = $newMdName
    } elseif ($addr -eq '$D$2') {
        $hl.TextToDisplay = $newDeXlfName
    }
}
